$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''46.660.36'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  +3.62%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''2.270.71'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  +0.01%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = '''  +0.12%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''302.52'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  +0.08%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''100.28'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  +6.30%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = '''  -0.44%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = '''1.00'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '''  +0.00%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = '''0.516'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''  +0.93%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = '''35.78'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''  +4.30%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = '''0.0784'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''  -0.74%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = '''  -0.35%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = '''  -0.90%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = '''2.622.33'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  +0.32%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = '''2.278.94'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  +0.47%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = '''13.68'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  -0.68%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = '''  +0.03%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = '''46.652.19'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''  +3.94%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = '''13.18'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''  +2.61%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = '''0.0₃0929'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''  +0.43%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = '''5.92'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  -3.09%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = '''65.32'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''  -0.29%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = '''247.29'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''  +3.20%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = '''  -1.51%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = '''0.998'
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Value = '''1.88'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''  -1.83%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = '''42.71'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '''  +2.31%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = '''2.24'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '''  -0.18%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = '''9.74'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '''  +1.91%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = '''19.88'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '''  +1.62%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = '''2.78'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '''  +7.92%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = '''  -1.86%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = '''147.26'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '''  -3.28%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = '''3.27'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '''  +10.72%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = '''0.0777'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '''  -1.50%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = '''  +10.09%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = '''  -0.70%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = '''  +15.84%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = '''1.73'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''  -1.66%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = '''  +1.56%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = '''3.26'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''  -0.13%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = '''  -3.28%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = '''  -0.01%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = '''  +0.53%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = '''1.817.90'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''  +2.81%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = '''90.29'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  +19.39%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = '''0.190'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''  -2.47%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = '''72.44'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''  +3.09%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = '''  +4.51%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = '''94.62'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''  -2.48%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = '''2.494.95'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '''  +0.06%  '
$ws.Range("E51").Style = "Normal"
